$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate current row 2 (values + styles) down into new rows 3 and 4,
# so the new rows inherit the same per-column styles (date/text/center).
$ws.Range("A2:E2").Copy($ws.Range("A3"))
$ws.Range("A2:E2").Copy($ws.Range("A4"))

# Row 2 becomes the new "OS Tick Initial Release" entry.
# Cell write order matters for shared-string interning order, so set
# C2/D2/B2 before A2/E2 to reproduce the expected sharedStrings layout.
$ws.Range("C2").Value = "JMR"
$ws.Range("D2").Value = "OS Tick Initial Release"
$ws.Range("B2").Value = "00"
$ws.Range("A2").Value = 41681
$ws.Range("E2").Value = "Done"

# Row 3 keeps the original "Fix notification..." entry values (already
# copied from the old row 2 above, so nothing else to do there).

# Row 4 is the new "OS Task Initial Release" entry.
$ws.Range("A4").Value = 41699
$ws.Range("B4").Value = "02"
$ws.Range("C4").Value = "SPA"
$ws.Range("D4").Value = "OS Task Initial Release"
$ws.Range("E4").Value = "On Process"

# Re-apply the text format to column B; this also promotes the header
# cell B1 to the new centered/wrapped text style.
$ws.Columns.Item(2).NumberFormat = "@"

# Move the active selection, matching the edited workbook's cursor spot.
$ws.Range("D6").Select()
